# Updates per "chore: update Sheets via scheduled runner"
# Each sheet's table columns are: A Leve Name, B Leve Item, C Leve Level,
# D Leve EXP, E Leve Gil, F Leve Amount, G Leve Item ID, H currentAveragePrice,
# I currentAveragePriceNQ, J currentAveragePriceHQ, K LevePriceNQ,
# L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
#
# NOTE: this runtime's Range.Value setter does not reliably apply array
# writes across multi-cell ranges, so every cell is written individually.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, $row, $h, $i, $j, $k, $l, $m, $n)
    $ws.Range("H$row").Value = $h
    $ws.Range("I$row").Value = $i
    $ws.Range("J$row").Value = $j
    $ws.Range("K$row").Value = $k
    $ws.Range("L$row").Value = $l
    $ws.Range("M$row").Value = $m
    $ws.Range("N$row").Value = $n
}

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

Set-Row $ws 18 6954620 13889238 20002 13889238 20002 -13888954 -20570

Set-Row $ws 116 4450.8887 3467.2222 5434.5557 3467.2222 5434.5557 -25.22220000000016 -12318.5557

Set-Row $ws 129 3206478 17857844 1491.875 53573532 4475.625 -53568532 -14475.625

Set-Row $ws 137 2275297.8 2705151.5 3214.2856 8115454.5 9642.856800000001 -8112904.5 -14742.8568

Set-Row $ws 138 2152.5823 1165.68 3854.138 3497.04 11562.414 1642.96 -21842.414

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

Set-Row $ws 2 14707416 16667658 5600 16667658 5600 -16667545 -5826

Set-Row $ws 116 14707416 16667658 5600 16667658 5600 -16665364 -10188

Set-Row $ws 122 3129.7827 2343.2727 3850.75 7029.8181 11552.25 -4579.8181 -16452.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

Set-Row $ws 3 14707416 16667658 5600 16667658 5600 -16667544 -5828

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

Set-Row $ws 10 6266.6924 549.9 25322.666 549.9 25322.666 -410.9 -25600.666

$ws.Range("H19").Value = 681.875
$ws.Range("I19").Value = 331
$ws.Range("K19").Value = 331
$ws.Range("M19").Value = -161

$ws.Range("H24").Value = 681.875
$ws.Range("I24").Value = 331
$ws.Range("K24").Value = 331
$ws.Range("M24").Value = -161

Set-Row $ws 31 2504469.8 3228961.8 8997.777 3228961.8 8997.777 -3228666.8 -9587.777

Set-Row $ws 34 2504469.8 3228961.8 8997.777 3228961.8 8997.777 -3228759.8 -9401.777

Set-Row $ws 132 2607.6155 1832.3334 4352 5497.0002 13056 -2967.0002 -18116

Set-Row $ws 134 1204.7593 818.2414 1653.12 2454.7242 4959.36 80.27579999999989 -10029.36

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

Set-Row $ws 5 591.78125 258.5926 2391 775.7778000000001 7173 -663.7778000000001 -7397

Set-Row $ws 135 591.78125 258.5926 2391 2327.3334 21519 207.6666 -26589

Set-Row $ws 137 2561.65 2140 2983.3 6420 8949.900000000001 -1320 -19149.9

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

Set-Row $ws 46 1485 436.77777 2533.2222 436.77777 2533.2222 -248.77777 -2909.2222

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H13").Value = 85006
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

Set-Row $ws 113 936.75 326.5 2462.375 979.5 7387.125 1190.5 -11727.125

Set-Row $ws 126 5558084 2023.3846 20003842 6070.1538 60011526 -3600.1538 -60016466

Set-Row $ws 132 196203.28 223510.53 59667 670531.59 179001 -668001.59 -184061

Write-Host "Applied Sheets update"
